$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.184.72"
$ws.Range("E2").Value = "'  -3.98%  "

$ws.Range("D3").Value = "'1.657.74"
$ws.Range("E3").Value = "'  -2.79%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "'  +0.23%  "

$ws.Range("D5").Value = "'217.06"
$ws.Range("E5").Value = "'  -3.06%  "

$ws.Range("D6").Value = "'0.5163"
$ws.Range("E6").Value = "'  -2.95%  "

$ws.Range("E7").Value = "'  +0.36%  "

$ws.Range("D8").Value = "'0.2581"
$ws.Range("E8").Value = "'  -3.02%  "

$ws.Range("D9").Value = "'0.06433"
$ws.Range("E9").Value = "'  -2.50%  "

$ws.Range("D10").Value = "'19.99"
$ws.Range("E10").Value = "'  -3.54%  "

$ws.Range("D11").Value = "'0.07786"
$ws.Range("E11").Value = "'  +1.61%  "

$ws.Range("D12").Value = "'1.657.13"
$ws.Range("E12").Value = "'  -3.00%  "

$ws.Range("D13").Value = "'4.296"
$ws.Range("E13").Value = "'  -4.61%  "

$ws.Range("D14").Value = "'1.886.17"
$ws.Range("E14").Value = "'  -2.78%  "

$ws.Range("D15").Value = "'0.5538"
$ws.Range("E15").Value = "'  -4.74%  "

$ws.Range("D16").Value = "'0.0₅8044"
$ws.Range("E16").Value = "'  -1.54%  "

$ws.Range("D17").Value = "'64.29"
$ws.Range("E17").Value = "'  -4.89%  "

$ws.Range("D18").Value = "'26.220.56"
$ws.Range("E18").Value = "'  -4.05%  "

$ws.Range("D19").Value = "'1.006"
$ws.Range("E19").Value = "'  +0.26%  "

$ws.Range("D20").Value = "'210.58"
$ws.Range("E20").Value = "'  -2.15%  "

$ws.Range("D21").Value = "'4.417"
$ws.Range("E21").Value = "'  -4.47%  "

$ws.Range("E22").Value = "'  -3.22%  "

$ws.Range("D23").Value = "'6.029"
$ws.Range("E23").Value = "'  +0.78%  "

$ws.Range("E24").Value = "'  +0.30%  "

$ws.Range("D25").Value = "'145.53"
$ws.Range("E25").Value = "'  +0.91%  "

$ws.Range("D26").Value = "'1.752"
$ws.Range("E26").Value = "'  +3.72%  "

$ws.Range("D27").Value = "'0.1170"
$ws.Range("E27").Value = "'  -2.76%  "

$ws.Range("D28").Value = "'6.986"
$ws.Range("E28").Value = "'  -3.27%  "

$ws.Range("E29").Value = "'  -2.54%  "

$ws.Range("D30").Value = "'0.05189"
$ws.Range("E30").Value = "'  -3.19%  "

$ws.Range("E31").Value = "'  -2.89%  "

$ws.Range("D32").Value = "'3.351"
$ws.Range("E32").Value = "'  -3.77%  "

$ws.Range("D33").Value = "'3.231"
$ws.Range("E33").Value = "'  -5.30%  "

$ws.Range("D34").Value = "'1.571"
$ws.Range("E34").Value = "'  -4.28%  "

$ws.Range("D35").Value = "'2.759"
$ws.Range("E35").Value = "'  -3.72%  "

$ws.Range("B36").Value = "'ARBITRUM"
$ws.Range("C36").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9318"
$ws.Range("E36").Value = "'  -1.83%  "

$ws.Range("B37").Value = "'HuobiToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.373"
$ws.Range("E37").Value = "'  -1.02%  "

$ws.Range("D38").Value = "'1.172.55"
$ws.Range("E38").Value = "'  +12.42%  "

$ws.Range("D39").Value = "'0.5700"
$ws.Range("E39").Value = "'  -2.38%  "

$ws.Range("D40").Value = "'0.01597"
$ws.Range("E40").Value = "'  -2.64%  "

$ws.Range("E41").Value = "'  +0.26%  "

$ws.Range("D42").Value = "'0.8402"
$ws.Range("E42").Value = "'  -0.24%  "

$ws.Range("D43").Value = "'5.674"
$ws.Range("E43").Value = "'  -2.37%  "

$ws.Range("D44").Value = "'100.41"
$ws.Range("E44").Value = "'  -0.46%  "

$ws.Range("D45").Value = "'1.796.51"
$ws.Range("E45").Value = "'  -2.80%  "

$ws.Range("E46").Value = "'  +1.04%  "

$ws.Range("D47").Value = "'0.4541"
$ws.Range("E47").Value = "'  +0.41%  "

$ws.Range("D48").Value = "'55.86"
$ws.Range("E48").Value = "'  -3.37%  "

$ws.Range("D49").Value = "'1.003"
$ws.Range("E49").Value = "'  -0.18%  "

$ws.Range("D50").Value = "'7.896"
$ws.Range("E50").Value = "'  -2.54%  "

$ws.Range("D51").Value = "'0.05064"
$ws.Range("E51").Value = "'  -3.27%  "
